$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7 ("Initial QuadSolver and references") gains an actual time ("30 min")
# and 100% complete, mirroring the pattern of the other already-completed rows.
$ws.Range("E7").Value = "30 min"
$ws.Range("F7").Value = 1
$ws.Range("F7").NumberFormat = "0%"

# Move the active selection to D17, matching the saved cursor position.
$ws.Range("D17").Select()
